$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to match the edited data
$ws.Range("C2").Value = "Rumenko"
$ws.Range("A3").Value = 421
$ws.Range("A4").Value = 201
$ws.Range("A5").Value = 421
$ws.Range("A6").Value = 342

# Update the active cell selection to A6
$ws.Range("A6").Select()
